$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the casing of three class names in column A
$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A4").Value = "mdaTitle"
$ws.Range("A8").Value = "pageTitleNewTab"

# Update the selected/active cell on the sheet
$ws.Range("A8").Select()
